# Competition.xlsx update — "regn, trading updates, crypto updates"
#
# 1) Main sheet: append row 20 (2022-08-08 snapshot) continuing the
#    running performance table, filling the D:G ratio/diff formulas
#    down from row 19.
# 2) Trades sheet: append the closing trade blocks for RBLX and SGHC
#    (entry leg + "Total X" roll-up row each) plus a combined "Total"
#    row, mirroring the existing blocks above them.
# 3) View state: Main's selection lands on the new B20 cell, then the
#    Trades sheet becomes the active tab with G293:L293 selected (the
#    new combined-total row).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Main")
$ws2 = $wb.Worksheets.Item("Trades")

# ---------------------------------------------------------------
# Main sheet: new row 20
# ---------------------------------------------------------------

# Carry the B:G formatting/styles down from row 19 onto row 20 first …
$ws1.Range("B19:G19").Copy($ws1.Range("B20:G20"))

# … then overwrite with row 20's own literal values / formulas.
$ws1.Range("B20").Value = 44781
$ws1.Range("C20").Value = 45951.47
$ws1.Range("D20").Formula = "=C20-C19"
$ws1.Range("E20").Formula = "=C20-`$C`$4"
$ws1.Range("F20").Formula = "=+C20/C19-1"
$ws1.Range("G20").Formula = "=C20/`$C`$4-1"

$ws1.Range("B20").Select()

# ---------------------------------------------------------------
# Trades sheet: RBLX close-out (rows 289-290)
# ---------------------------------------------------------------

$ws2.Range("E289:K291").NumberFormat = "#,##0.00"

$ws2.Range("B289").Value = "RBLX"
$ws2.Range("C289").Value = "2022-08-08, 10:20:07"
$ws2.Range("D289").Value = -50
$ws2.Range("E289").Value = 49.88
$ws2.Range("F289").Value = 48.9
$ws2.Range("G289").Value = 2494
$ws2.Range("H289").Value = -1.06
$ws2.Range("I289").Value = -2030.5
$ws2.Range("J289").Value = 462.44
$ws2.Range("K289").Value = 49
$ws2.Range("L289").Value = "C;P"

$ws2.Range("B290").Value = "Total RBLX"
$ws2.Range("D290").Value = -50
$ws2.Range("E290").Value = " "
$ws2.Range("G290").Value = 2494
$ws2.Range("H290").Value = -1.06
$ws2.Range("I290").Value = -2030.5
$ws2.Range("J290").Value = 462.44
$ws2.Range("K290").Value = 49
$ws2.Range("L290").NumberFormat = "#,##0.00"
$ws2.Range("L290").Value = " "

# ---------------------------------------------------------------
# Trades sheet: SGHC close-out (rows 291-292)
# ---------------------------------------------------------------

$ws2.Range("B291").Value = "SGHC"
$ws2.Range("C291").Value = "2022-08-08, 09:36:37"
$ws2.Range("D291").Value = 500
$ws2.Range("E291").Value = 5.3196000000000003
$ws2.Range("F291").Value = 5.35
$ws2.Range("G291").Value = -2659.8
$ws2.Range("H291").Value = -2.5
$ws2.Range("I291").Value = 2047.39
$ws2.Range("J291").Value = -614.91
$ws2.Range("K291").Value = 15.2
$ws2.Range("L291").Value = "C;P"

$ws2.Range("B292").Value = "Total SGHC"
$ws2.Range("D292").Value = 500
$ws2.Range("E292").NumberFormat = "#,##0.00"
$ws2.Range("E292").Value = " "
$ws2.Range("G292").NumberFormat = "#,##0.00"
$ws2.Range("G292").Value = -2659.8
$ws2.Range("H292").NumberFormat = "#,##0.00"
$ws2.Range("H292").Value = -2.5
$ws2.Range("I292").NumberFormat = "#,##0.00"
$ws2.Range("I292").Value = 2047.39
$ws2.Range("J292").NumberFormat = "#,##0.00"
$ws2.Range("J292").Value = -614.91
$ws2.Range("K292").NumberFormat = "#,##0.00"
$ws2.Range("K292").Value = 15.2
$ws2.Range("L292").NumberFormat = "#,##0.00"
$ws2.Range("L292").Value = " "

# ---------------------------------------------------------------
# Trades sheet: combined grand-total row (row 293)
# ---------------------------------------------------------------

$ws2.Range("B293").Value = "Total"
$ws2.Range("G293").Value = -165.8
$ws2.Range("H293").Value = -3.56
$ws2.Range("I293:L293").NumberFormat = "#,##0.00"
$ws2.Range("I293").Value = 16.89
$ws2.Range("J293").Value = -152.47999999999999
$ws2.Range("K293").Value = 64.2
$ws2.Range("L293").Value = " "

# ---------------------------------------------------------------
# Final view state: Trades becomes the active tab, new total row
# selected.
# ---------------------------------------------------------------

$ws2.Activate()
$ws2.Range("G293:L293").Select()
